# Add a "Hungary" market tab, modelled on the existing "Slovakia" tab
# (same layout/styles, just a new market name + part number), matching
# the "Added HungaryFC Test data" commit.

$wb = $excel.ActiveWorkbook

# The Slovakia sheet is the last tab and the template for every other
# per-country sheet in this workbook.
$slovakia = $wb.Sheets.Item("Slovakia")

# Move-or-Copy "(create a copy)" placed immediately after Slovakia -
# this duplicates all rows/styles/merged cells/page setup.
$slovakia.Copy([System.Reflection.Missing]::Value, $slovakia)
$hungary = $wb.Sheets.Item($slovakia.Index + 1)
$hungary.Name = "Hungary"

# Update the two market-specific cells on the new sheet.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3599"

# Restore Slovakia's own selection state (copying the sheet leaves the
# whole sheet selected on the source tab) and move the active
# selection/tab over to the newly created Hungary sheet.
$slovakia.Activate()
$slovakia.Range("A1:XFD1048576").Select() | Out-Null

$hungary.Activate()
$hungary.Range("B2:B4").Select() | Out-Null
$excel.ActiveWindow.RangeSelection.Activate() | Out-Null
